# movements.docx template edit:
#  1. First "التاريخ: $DATE_AR$" placeholder (table header cell) -> $TOD_AR$
#  2. $LOGO$ paragraph (table header cell) gets forced LTR reading order
#     (<w:bidi w:val="0"/> added to its pPr, right after pStyle)
#  3. Title run "يومية تحركات عن يوم (" is split so the hard-coded
#     "تحركات" becomes a $DOC_TYP$ placeholder run, keeping the
#     surrounding Arabic text (and its rtl/hint=cs formatting) intact.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) $DATE_AR$ -> $TOD_AR$ in the "التاريخ:" line only (paragraph 6).
#    The title further down (paragraph 13) also contains a $DATE_AR$
#    placeholder that must stay untouched, so the Find/Replace is scoped
#    to that single paragraph's Range and uses wdReplaceOne (1), not
#    wdReplaceAll, to avoid spilling outside of it.
# ---------------------------------------------------------------------
$dateParagraph = $d.Paragraphs.Item(6)
$dateRange = $dateParagraph.Range
$dateRange.Find.ClearFormatting()
$dateRange.Find.Execute("`$DATE_AR`$", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "`$TOD_AR`$", 1)

# ---------------------------------------------------------------------
# 2) Force LTR reading order on the $LOGO$ paragraph. Plain
#    ParagraphFormat.ReadingOrder assignment does not persist for
#    paragraphs that live inside a table cell in this runtime, so the
#    paragraph is rewritten in place via InsertXML with the single
#    <w:bidi w:val="0"/> addition applied (schema-ordered right after
#    pStyle, before jc) and nothing else changed.
# ---------------------------------------------------------------------
$logoParagraph = $d.Paragraphs.Item(9)
$logoRange = $logoParagraph.Range
$logoXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0BD2A4B3" w14:textId="77777777" w:rsidR="008F7506" w:rsidRDefault="008F7506" w:rsidP="00866FB9"><w:pPr><w:pStyle w:val="Header"/><w:bidi w:val="0"/><w:jc w:val="center"/><w:rPr><w:rtl/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:noProof/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>$LOGO$</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$logoRange.InsertXML($logoXml)

# ---------------------------------------------------------------------
# 3) Split the title's first run ("يومية تحركات عن يوم (") into three
#    runs so the document-type word becomes a $DOC_TYP$ placeholder:
#      "يومية "  (unchanged rtl/hint=cs formatting)
#      "$DOC_TYP$" (new run, LTR placeholder formatting, bold+underline)
#      " عن يوم (" (unchanged rtl/hint=cs formatting)
#    The remaining runs in the paragraph ($WEEKDAY$, ") الموافق ",
#    $DATE_AR$) are re-emitted unchanged.
# ---------------------------------------------------------------------
$titleParagraph = $d.Paragraphs.Item(13)
$titleRange = $titleParagraph.Range
$titleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="3FD7B792" w14:textId="28BF3759" w:rsidR="006C4641" w:rsidRPr="00145020" w:rsidRDefault="00D20AF4" w:rsidP="004D0E50"><w:pPr><w:spacing w:after="0"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Simplified Arabic" w:hAnsi="Simplified Arabic" w:cs="Simplified Arabic"/><w:b/><w:bCs/><w:sz w:val="12"/><w:szCs w:val="12"/><w:u w:val="single"/><w:rtl/><w:lang w:bidi="ar-EG"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Simplified Arabic" w:hAnsi="Simplified Arabic" w:cs="Simplified Arabic" w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/><w:rtl/><w:lang w:bidi="ar-EG"/></w:rPr><w:t xml:space="preserve">يومية </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Simplified Arabic" w:hAnsi="Simplified Arabic" w:cs="Simplified Arabic"/><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/><w:lang w:bidi="ar-EG"/></w:rPr><w:t>$DOC_TYP$</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Simplified Arabic" w:hAnsi="Simplified Arabic" w:cs="Simplified Arabic" w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/><w:rtl/><w:lang w:bidi="ar-EG"/></w:rPr><w:t xml:space="preserve"> عن يوم (</w:t></w:r><w:r w:rsidR="00E15714"><w:rPr><w:rFonts w:ascii="Simplified Arabic" w:hAnsi="Simplified Arabic" w:cs="Simplified Arabic"/><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/><w:lang w:bidi="ar-EG"/></w:rPr><w:t>$WEEKDAY$</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Simplified Arabic" w:hAnsi="Simplified Arabic" w:cs="Simplified Arabic" w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/><w:rtl/><w:lang w:bidi="ar-EG"/></w:rPr><w:t xml:space="preserve">) الموافق </w:t></w:r><w:r w:rsidR="00E15714"><w:rPr><w:rFonts w:ascii="Simplified Arabic" w:hAnsi="Simplified Arabic" w:cs="Simplified Arabic"/><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/><w:lang w:bidi="ar-EG"/></w:rPr><w:t>$DATE_AR$</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$titleRange.InsertXML($titleXml)

Write-Output "movements template updated"
